# CRE21-020: To include the School Code in the Pre-Authorization Checking File
# (for both PPP-PS and PPP-KG)
#
# 1) Sheet "03": the header formula in K3 must now show "School Code" when
#    there is no RCH Type (L4) but there is a value in K4, and keep showing
#    "RCH Code" when L4 is also populated.
# 2) Sheet "Change History": log the change as a new row (Item 7 / CRE21-020).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the K3 header formula on the "03" worksheet
# ---------------------------------------------------------------------------
$ws03 = $wb.Worksheets.Item("03")
$ws03.Range("K3").Formula = '=IF(LEN(K4) > 0,IF(LEN(L4) > 0,"RCH Code","School Code"),"")'

# ---------------------------------------------------------------------------
# 2. Append a new entry to the "Change History" worksheet
# ---------------------------------------------------------------------------
$wsCH = $wb.Worksheets.Item("Change History")

# Base the new row on the previous entry (row 9) so that it inherits the same
# look & feel (fonts, borders, number format, wrap text, etc.) before we
# overwrite the actual values.
$wsCH.Range("A9:D9").Copy($wsCH.Range("A10:D10"))
$wsCH.Rows.Item(10).RowHeight = 30

$wsCH.Range("A10").Value = 7
$wsCH.Range("B10").Value = "CRE21-020"
$wsCH.Range("C10").Value = "To include the School Code in the Pre-Authorization Checking File (for both PPP-PS and PPP-KG)"
$wsCH.Range("D10").Value = 44530
